$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.866.83"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.39%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.813.37"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.96%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.19%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("E5").Value = "'  +0.23%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +0.23%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.4648"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +1.66%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.3685"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -0.85%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.07366"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +1.72%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.8690"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +1.64%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  +0.01%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'1.814.53"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +1.00%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'5.345"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +0.76%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.07072"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +0.48%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'91.76"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +1.64%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'6.491"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.10%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'1.002"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +0.17%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'0.000008685"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +0.63%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'1.001"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +0.23%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'14.75"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +0.88%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'26.897.11"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +0.46%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'5.330"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +0.83%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'10.55"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -0.65%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'2.041.34"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +1.31%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  -0.37%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'151.09"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +1.14%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'2.174"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +0.63%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  +0.70%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  +2.15%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'115.47"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +1.15%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'0.08925"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +1.01%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'0.7654"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +1.07%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'1.160"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +0.07%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'4.501"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +1.39%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'2.903"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +0.77%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'1.001"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +0.23%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  -2.26%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  +0.94%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.05283"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +1.31%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'2.944"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +1.48%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'7.264"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +1.87%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.5309"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +1.28%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'2.353"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -0.78%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  +1.00%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'8.410"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -0.84%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.4923"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -3.00%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'10.47"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +2.50%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  +0.32%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").Value = "'1.667"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +1.21%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D50").Value = "'103.56"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -0.66%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.06274"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -0.32%  "
$ws.Range("E51").Style = "Normal"
